# Marksheet update: correct the "Marking"/"Total" row totals and the
# correct/total marks fraction shown in the Max column of the Total row.
$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "quiz") {
        $ws = $sheet
        break
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

# Row 11 "Marking" total correct count: 3 -> 5
$ws.Range("B11").Value = 5

# Row 12 "Total" correct count: 33 -> 55
$ws.Range("B12").Value = 55

# Row 12 "Total" correct/total marks fraction: 27/84 -> 55/140
$ws.Range("E12").Value = "55/140"
